$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A2: shift the starting timestamp forward (DST-threshold refresh) ---
$ws.Range("A2").Value = 44501.291666666664

# --- Pre-seed the date format/style on the brand-new rows (17-24) by
#     copying the existing date cell's format so the later Formula
#     assignment doesn't fall back to "General" on newly-created cells ---
$ws.Range("A16").Copy()
$ws.Range("A17:A24").PasteSpecial(-4122)

# --- A5: continue the existing (non-shared) per-cell MROUND formula ---
$ws.Range("A5").Formula = '=MROUND(A4+1, "01:00")'

# --- A6:A24: assign as one multi-cell formula -> becomes an Excel shared formula group ---
$ws.Range("A6:A24").Formula = '=MROUND(A5+1, "01:00")'

# --- New data rows: B5:C19 ---
$ws.Range("B5").Value = 400
$ws.Range("C5").Value = 40
$ws.Range("B6").Value = 500
$ws.Range("C6").Value = 50
$ws.Range("B7").Value = 600
$ws.Range("C7").Value = 60
$ws.Range("B8").Value = 700
$ws.Range("C8").Value = 70
$ws.Range("B9").Value = 800
$ws.Range("C9").Value = 80
$ws.Range("B10").Value = 900
$ws.Range("C10").Value = 90
$ws.Range("B11").Value = 1000
$ws.Range("C11").Value = 100
$ws.Range("B12").Value = 1100
$ws.Range("C12").Value = 110
$ws.Range("B13").Value = 1200
$ws.Range("C13").Value = 120
$ws.Range("B14").Value = 1300
$ws.Range("C14").Value = 130
$ws.Range("B15").Value = 1400
$ws.Range("C15").Value = 140
$ws.Range("B16").Value = 1500
$ws.Range("C16").Value = 150
$ws.Range("B17").Value = 1600
$ws.Range("C17").Value = 160
$ws.Range("B18").Value = 1700
$ws.Range("C18").Value = 170
$ws.Range("B19").Value = 1800
$ws.Range("C19").Value = 180

# --- Move the active selection to B25 (next empty row in column B) ---
$null = $ws.Range("B25").Select()
